# Updates cryptocurrency price and volume(1h) values in the worksheet
# to reflect the latest scraped data (GitHub Actions scheduled update).
#
# Numeric-looking price strings are written via NumberFormat "@" (text)
# so Excel keeps them as text (matching the source data) instead of
# auto-converting them to numbers; the style is then reset back to the
# workbook's default "Normal" style so no new per-cell formatting is
# introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "46.612.26"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.272.67"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "300.56"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "100.31"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.87%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.574"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.19%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.510"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.40%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "35.11"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.20%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0800"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.77%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.09"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "2.618.61"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "2.271.83"
$ws.Range("E15").Value = "  -2.11%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "13.67"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").Value = "46.641.56"
$ws.Range("E17").Value = "  +1.65%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.797"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.06%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.70"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("E21").Value = "  -5.52%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "66.20"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "248.25"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -5.72%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "41.62"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("E28").Value = "  -1.67%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.67"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "20.24"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  +7.51%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.36"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +11.17%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "147.16"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("E36").Value = "  +7.46%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.116"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  +12.36%  "
$ws.Range("E39").Value = "  -6.51%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.87"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.17%  "
$ws.Range("E41").Value = "  -6.30%  "
$ws.Range("E42").Value = "  -6.62%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "92.79"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +16.41%  "
$ws.Range("D45").Value = "1.787.03"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  -5.97%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "71.21"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("E48").Value = "  -6.19%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.81"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "94.92"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.08%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.88"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.41%  "
